$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-10-13 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-10-14 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("34×78=2652", $true, $false, $false, $false, $false, $true, 1, $false, "52×44=2288", 2) | Out-Null
$d.Content.Find.Execute("50×72=3600", $true, $false, $false, $false, $false, $true, 1, $false, "94×96=9024", 2) | Out-Null
$d.Content.Find.Execute("33×41=1353", $true, $false, $false, $false, $false, $true, 1, $false, "87×39=3393", 2) | Out-Null
$d.Content.Find.Execute("62×83=5146", $true, $false, $false, $false, $false, $true, 1, $false, "88×16=1408", 2) | Out-Null
$d.Content.Find.Execute("84×89=7476", $true, $false, $false, $false, $false, $true, 1, $false, "75×13=975", 2) | Out-Null
$d.Content.Find.Execute("30×79=2370", $true, $false, $false, $false, $false, $true, 1, $false, "74×78=5772", 2) | Out-Null
$d.Content.Find.Execute("85×89=7565", $true, $false, $false, $false, $false, $true, 1, $false, "54×73=3942", 2) | Out-Null
$d.Content.Find.Execute("81×75=6075", $true, $false, $false, $false, $false, $true, 1, $false, "57×78=4446", 2) | Out-Null
$d.Content.Find.Execute("21×30=630", $true, $false, $false, $false, $false, $true, 1, $false, "15×50=750", 2) | Out-Null
$d.Content.Find.Execute("44×19=836", $true, $false, $false, $false, $false, $true, 1, $false, "26×97=2522", 2) | Out-Null
$d.Content.Find.Execute("87×93=8091", $true, $false, $false, $false, $false, $true, 1, $false, "16×57=912", 2) | Out-Null
$d.Content.Find.Execute("15×22=330", $true, $false, $false, $false, $false, $true, 1, $false, "94×93=8742", 2) | Out-Null
$d.Content.Find.Execute("57×14=798", $true, $false, $false, $false, $false, $true, 1, $false, "59×17=1003", 2) | Out-Null
$d.Content.Find.Execute("16×28=448", $true, $false, $false, $false, $false, $true, 1, $false, "21×12=252", 2) | Out-Null
$d.Content.Find.Execute("50×25=1250", $true, $false, $false, $false, $false, $true, 1, $false, "25×34=850", 2) | Out-Null
$d.Content.Find.Execute("57×72=4104", $true, $false, $false, $false, $false, $true, 1, $false, "28×60=1680", 2) | Out-Null
$d.Content.Find.Execute("86×22=1892", $true, $false, $false, $false, $false, $true, 1, $false, "93×81=7533", 2) | Out-Null
$d.Content.Find.Execute("82×95=7790", $true, $false, $false, $false, $false, $true, 1, $false, "54×20=1080", 2) | Out-Null
$d.Content.Find.Execute("25×84=2100", $true, $false, $false, $false, $false, $true, 1, $false, "49×19=931", 2) | Out-Null
$d.Content.Find.Execute("32×36=1152", $true, $false, $false, $false, $false, $true, 1, $false, "19×83=1577", 2) | Out-Null
$d.Content.Find.Execute("33×99=3267", $true, $false, $false, $false, $false, $true, 1, $false, "25×85=2125", 2) | Out-Null
$d.Content.Find.Execute("98×65=6370", $true, $false, $false, $false, $false, $true, 1, $false, "19×17=323", 2) | Out-Null
$d.Content.Find.Execute("22×19=418", $true, $false, $false, $false, $false, $true, 1, $false, "44×48=2112", 2) | Out-Null
$d.Content.Find.Execute("88×56=4928", $true, $false, $false, $false, $false, $true, 1, $false, "61×79=4819", 2) | Out-Null
$d.Content.Find.Execute("91×97=8827", $true, $false, $false, $false, $false, $true, 1, $false, "82×97=7954", 2) | Out-Null
